# Updates the "Estado de Cuenta" data table (rows 16-86) on Hoja1:
#  - Re-sorts the records by Periodo Mora (ascending) instead of by worker
#    then Periodo Mora (descending); within the same period, ADRIANA's row
#    comes before ANA SUSANA's row (she has no record for period 1704).
#  - Updates "Salario Basico" (column G) from 737717 to 781242 for every row.
# "N Doc Trabajador" (C) and "Periodo Mora" (E) are text-formatted columns,
# so the values are assigned as strings to avoid Excel re-typing them as
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(16, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1704", 29509, 781242),
    @(17, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1705", 29509, 781242),
    @(18, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1705", 29509, 781242),
    @(19, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1706", 29509, 781242),
    @(20, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1706", 29509, 781242),
    @(21, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1707", 29509, 781242),
    @(22, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1707", 29509, 781242),
    @(23, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1708", 29509, 781242),
    @(24, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1708", 29509, 781242),
    @(25, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1709", 29509, 781242),
    @(26, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1709", 29509, 781242),
    @(27, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1710", 29509, 781242),
    @(28, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1710", 29509, 781242),
    @(29, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1711", 29509, 781242),
    @(30, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1711", 29509, 781242),
    @(31, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1712", 29509, 781242),
    @(32, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1712", 29509, 781242),
    @(33, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1801", 29509, 781242),
    @(34, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1801", 29509, 781242),
    @(35, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1802", 29509, 781242),
    @(36, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1802", 29509, 781242),
    @(37, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1803", 29509, 781242),
    @(38, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1803", 29509, 781242),
    @(39, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1804", 29509, 781242),
    @(40, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1804", 29509, 781242),
    @(41, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1805", 29509, 781242),
    @(42, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1805", 29509, 781242),
    @(43, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1806", 29509, 781242),
    @(44, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1806", 29509, 781242),
    @(45, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1807", 29509, 781242),
    @(46, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1807", 29509, 781242),
    @(47, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1808", 29509, 781242),
    @(48, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1808", 29509, 781242),
    @(49, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1809", 31249, 781242),
    @(50, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1809", 31249, 781242),
    @(51, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1810", 31249, 781242),
    @(52, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1810", 31249, 781242),
    @(53, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1811", 31249, 781242),
    @(54, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1811", 31249, 781242),
    @(55, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1812", 31249, 781242),
    @(56, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1812", 31249, 781242),
    @(57, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1901", 31249, 781242),
    @(58, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1901", 31249, 781242),
    @(59, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1902", 31249, 781242),
    @(60, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1902", 31249, 781242),
    @(61, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1903", 31249, 781242),
    @(62, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1903", 31249, 781242),
    @(63, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1904", 31249, 781242),
    @(64, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1904", 31249, 781242),
    @(65, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1905", 31249, 781242),
    @(66, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1905", 31249, 781242),
    @(67, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1906", 31249, 781242),
    @(68, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1906", 31249, 781242),
    @(69, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1907", 31249, 781242),
    @(70, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1907", 31249, 781242),
    @(71, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1908", 31249, 781242),
    @(72, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1908", 31249, 781242),
    @(73, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1909", 31249, 781242),
    @(74, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1909", 31249, 781242),
    @(75, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1910", 31249, 781242),
    @(76, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1910", 31249, 781242),
    @(77, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1911", 31249, 781242),
    @(78, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1911", 31249, 781242),
    @(79, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "1912", 31249, 781242),
    @(80, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "1912", 31249, 781242),
    @(81, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "2001", 31249, 781242),
    @(82, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "2001", 31249, 781242),
    @(83, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "2002", 31249, 781242),
    @(84, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "2002", 31249, 781242),
    @(85, "CC", "1047371667", "ADRIANA MARCELA MURRA FALLA", "2003", 31249, 781242),
    @(86, "CC", "1128047638", "ANA SUSANA DE VALENCIA SPATH", "2003", 31249, 781242)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $row[2]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[3]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[4]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[5]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $row[6]   # G: Salario Basico
}
